$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.264.08'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.014.74'
$ws.Range("E3").Value = '  +6.06%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.07'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.658'
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.18'
$ws.Range("E8").Value = '  +2.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.52'
$ws.Range("E9").Value = '  +5.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.358'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0711'
$ws.Range("E11").Value = '  -5.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0980'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.18'
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.307.79'
$ws.Range("E14").Value = '  +6.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.799'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.014.30'
$ws.Range("E16").Value = '  +5.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.85'
$ws.Range("E17").Value = '  -3.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.336.30'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.81'
$ws.Range("E19").Value = '  -4.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0806'
$ws.Range("E20").Value = '  -3.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '235.58'
$ws.Range("E21").Value = '  -4.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.64'
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.86'
$ws.Range("E23").Value = '  -6.50%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  -10.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.77'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.57'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.55'
$ws.Range("E28").Value = '  +6.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.93'
$ws.Range("E29").Value = '  -10.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.120'
$ws.Range("E30").Value = '  -6.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.42'
$ws.Range("E31").Value = '  +53.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.31'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0576'
$ws.Range("E33").Value = '  -5.16%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.86'
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0848'
$ws.Range("E36").Value = '  +15.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.94'
$ws.Range("E37").Value = '  -7.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  +7.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.850'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  -12.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0213'
$ws.Range("E41").Value = '  -6.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.96'
$ws.Range("E42").Value = '  -4.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.09'
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +16.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.61'
$ws.Range("E45").Value = '  -8.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.301.36'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.211.39'
$ws.Range("E49").Value = '  +6.41%  '
$ws.Range("E50").Value = '  -7.62%  '
$ws.Range("E51").Value = '  +14.92%  '
